$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.731.65"
$ws.Range("E2").Value = "  +1.18%  "
$ws.Range("D3").Value = "2.296.75"
$ws.Range("E3").Value = "  +0.93%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "497.29"
$ws.Range("E5").Value = "  +1.22%  "
$ws.Range("D6").Value = "128.39"
$ws.Range("E6").Value = "  +0.41%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").Value = "0.531"
$ws.Range("E8").Value = "  +1.10%  "
$ws.Range("D9").Value = "2.296.40"
$ws.Range("E9").Value = "  +0.82%  "
$ws.Range("D10").Value = "0.0951"
$ws.Range("E10").Value = "  +2.16%  "
$ws.Range("E11").Value = "  +2.47%  "
$ws.Range("D12").Value = "0.324"
$ws.Range("E12").Value = "  +2.42%  "
$ws.Range("D13").Value = "4.65"
$ws.Range("E13").Value = "  -2.69%  "
$ws.Range("D14").Value = "2.705.64"
$ws.Range("E14").Value = "  +0.99%  "
$ws.Range("D15").Value = "21.84"
$ws.Range("E15").Value = "  +3.07%  "
$ws.Range("D16").Value = "54.673.26"
$ws.Range("E16").Value = "  +1.14%  "
$ws.Range("E17").Value = "  +0.79%  "
$ws.Range("D18").Value = "2.305.27"
$ws.Range("E18").Value = "  +0.72%  "
$ws.Range("D19").Value = "10.10"
$ws.Range("E19").Value = "  +4.47%  "
$ws.Range("E20").Value = "  +3.26%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "6.49"
$ws.Range("E21").Value = "  +5.41%  "
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").Value = "302.76"
$ws.Range("E22").Value = "  -0.21%  "
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.47%  "
$ws.Range("E24").Value = "  -1.82%  "
$ws.Range("D25").Value = "62.88"
$ws.Range("E25").Value = "  -1.63%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("E27").Value = "  +2.03%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.407.98"
$ws.Range("E28").Value = "  +0.76%  "
$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").Value = "0.152"
$ws.Range("E29").Value = "  +4.80%  "
$ws.Range("D30").Value = "7.12"
$ws.Range("E30").Value = "  +0.65%  "
$ws.Range("D31").Value = "169.73"
$ws.Range("E31").Value = "  +0.39%  "
$ws.Range("D32").Value = "0.0₃0694"
$ws.Range("E32").Value = "  -0.77%  "
$ws.Range("D33").Value = "1.60"
$ws.Range("E33").Value = "  -0.35%  "
$ws.Range("D34").Value = "5.90"
$ws.Range("E34").Value = "  +1.91%  "
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").Value = "1.08"
$ws.Range("E37").Value = "  +0.65%  "
$ws.Range("D38").Value = "17.67"
$ws.Range("E38").Value = "  +0.62%  "
$ws.Range("E39").Value = "  +2.63%  "
$ws.Range("D40").Value = "0.871"
$ws.Range("E40").Value = "  +4.04%  "
$ws.Range("D41").Value = "3.69"
$ws.Range("E41").Value = "  +2.01%  "
$ws.Range("D42").Value = "35.55"
$ws.Range("E42").Value = "  -0.71%  "
$ws.Range("D43").Value = "1.41"
$ws.Range("E43").Value = "  +2.77%  "
$ws.Range("D44").Value = "0.376"
$ws.Range("E44").Value = "  +2.27%  "
$ws.Range("D45").Value = "3.36"
$ws.Range("E45").Value = "  +1.17%  "
$ws.Range("D46").Value = "128.22"
$ws.Range("E46").Value = "  +4.61%  "
$ws.Range("D47").Value = "4.84"
$ws.Range("E47").Value = "  +2.69%  "
$ws.Range("E48").Value = "  +1.28%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "0.549"
$ws.Range("E49").Value = "  +1.28%  "
$ws.Range("B50").Value = "Bittensor"
$ws.Range("C50").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D50").Value = "242.32"
$ws.Range("E50").Value = "  +1.09%  "
$ws.Range("D51").Value = "0.0486"
$ws.Range("E51").Value = "  +2.54%  "
